$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feedback")

$ws.Range("B2").Value = "Nirmal Kumar Pant"
$ws.Range("J2").Value = "Auto test ing 11"
$ws.Range("B3").Value = "Auto testing 8"
$ws.Range("E3").Value = "AutoTest 5we"

$ws.Range("E3").Select()
